# "all names and months in final excel alligned"
# Adds 11 new payroll rows (18-28) to the worksheet, reusing the existing
# "Pan yes" label already in row 18/19 (C18/C19, D18/D19 already contain the
# right values) and introducing four brand new labels ("Pan hes", "Pan kes",
# "pan qes", "pan bes") for the following rows. Column L gets a lohnart-code
# (shared string) tag on every one of the newly touched rows. Finally the
# sheet view is scrolled/zoomed and the selection is moved to the new bottom
# of the data block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18/19 already have C/D filled in ("Pan yes" / 8) - only L is new.
$ws.Range("L18").Value = "Netto"
$ws.Range("L19").Value = "Netto"

# Row 20-22: "Pan hes"
$ws.Range("C20").Value = "Pan hes"
$ws.Range("D20").Value = 4
$ws.Range("L20").Value = "GV pro Stunde"

$ws.Range("C21").Value = "Pan hes"
$ws.Range("D21").Value = 3
$ws.Range("L21").Value = "Netto"

$ws.Range("C22").Value = "Pan hes"
$ws.Range("D22").Value = 3
$ws.Range("L22").Value = "Netto"

# Row 23-24: "Pan kes"
$ws.Range("C23").Value = "Pan kes"
$ws.Range("D23").Value = 9
$ws.Range("L23").Value = "Netto"

$ws.Range("C24").Value = "Pan kes"
$ws.Range("D24").Value = 4
$ws.Range("L24").Value = "GV pro Stunde"

# Row 25-26: "pan qes"
$ws.Range("C25").Value = "pan qes"
$ws.Range("D25").Value = 6
$ws.Range("L25").Value = "Netto"

$ws.Range("C26").Value = "pan qes"
$ws.Range("D26").Value = 6
$ws.Range("L26").Value = "Netto"

# Row 27-28: "pan bes"
$ws.Range("C27").Value = "pan bes"
$ws.Range("D27").Value = 8
$ws.Range("L27").Value = "Netto"

$ws.Range("C28").Value = "pan bes"
$ws.Range("D28").Value = 8
$ws.Range("L28").Value = "Netto"

# Scroll/zoom the sheet view to frame the newly added rows and move the
# active selection down to the new bottom block (L25:L28, anchored at L25).
$aw = $ws.Application.ActiveWindow
$aw.ScrollColumn = 1
$aw.ScrollRow = 9
$aw.Zoom = 85

$ws.Range("L25:L28").Select()
